$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = "'26.452.24"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.58%  '
# Row 3
$ws.Range('D3').Value = "'1.676.07"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.54%  '
# Row 4
$ws.Range('E4').Value = '  -0.01%  '
# Row 5
$ws.Range('D5').Value = "'217.05"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.86%  '
# Row 6
$ws.Range('D6').Value = "'0.5317"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.53%  '
# Row 7
$ws.Range('E7').Value = '  -0.04%  '
# Row 8
$ws.Range('E8').Value = '  +3.81%  '
# Row 9
$ws.Range('D9').Value = "'0.06408"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.85%  '
# Row 10
$ws.Range('D10').Value = "'21.82"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.45%  '
# Row 11
$ws.Range('D11').Value = "'0.07804"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.83%  '
# Row 12
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = "'1.683.92"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.08%  '
# Row 13
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = "'4.514"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.33%  '
# Row 14
$ws.Range('D14').Value = "'0.5572"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.92%  '
# Row 15
$ws.Range('D15').Value = "'0.0₅8326"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.26%  '
# Row 16
$ws.Range('D16').Value = "'65.70"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.19%  '
# Row 17
$ws.Range('D17').Value = "'26.512.74"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.89%  '
# Row 18
$ws.Range('E18').Value = '  -0.06%  '
# Row 19
$ws.Range('D19').Value = "'4.780"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.93%  '
# Row 20
$ws.Range('D20').Value = "'193.68"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.07%  '
# Row 21
$ws.Range('E21').Value = '  +1.19%  '
# Row 22
$ws.Range('D22').Value = "'6.343"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.08%  '
# Row 23
$ws.Range('E23').Value = '  +0.00%  '
# Row 24
$ws.Range('D24').Value = "'142.40"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.69%  '
# Row 25
$ws.Range('D25').Value = "'0.1281"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +5.36%  '
# Row 26
$ws.Range('D26').Value = "'7.407"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.07%  '
# Row 27
$ws.Range('D27').Value = "'16.28"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.23%  '
# Row 28
$ws.Range('D28').Value = "'1.449"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.96%  '
# Row 29
$ws.Range('D29').Value = "'0.06290"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +5.76%  '
# Row 30
$ws.Range('D30').Value = "'1.273"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.68%  '
# Row 31
$ws.Range('D31').Value = "'3.625"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +5.44%  '
# Row 32
$ws.Range('D32').Value = "'3.453"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.28%  '
# Row 33
$ws.Range('D33').Value = "'1.681"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.70%  '
# Row 34
$ws.Range('D34').Value = "'1.009"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.63%  '
# Row 35
$ws.Range('D35').Value = "'0.6206"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +9.44%  '
# Row 36
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').Value = "'2.416"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.85%  '
# Row 37
$ws.Range('B37').Value = 'MXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D37').Value = "'2.786"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.92%  '
# Row 38
$ws.Range('D38').Value = "'6.184"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +7.68%  '
# Row 39
$ws.Range('D39').Value = "'0.01638"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.49%  '
# Row 40
$ws.Range('D40').Value = "'1.096.89"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +6.11%  '
# Row 41
$ws.Range('D41').Value = "'0.8658"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.97%  '
# Row 42
$ws.Range('E42').Value = '  -0.08%  '
# Row 43
$ws.Range('D43').Value = "'100.59"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.24%  '
# Row 44
$ws.Range('D44').Value = "'1.821.65"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.13%  '
# Row 45
$ws.Range('D45').Value = "'57.80"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.96%  '
# Row 46
$ws.Range('D46').Value = "'8.133"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.22%  '
# Row 47
$ws.Range('D47').Value = "'1.003"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.33%  '
# Row 48
$ws.Range('D48').Value = "'0.0₈104"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -6.26%  '
# Row 49
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').Value = "'0.05213"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.02%  '
# Row 50
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').Value = "'1.486"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +7.30%  '
# Row 51
$ws.Range('D51').Value = "'6.052"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.19%  '
